# Insert a new "Table of Contents" Heading2 paragraph immediately before
# the existing "[ToC]" placeholder paragraph.

$d = $word.ActiveDocument

# Locate the "[ToC]" placeholder text using Find (literal match, no
# wildcards) so the square brackets are treated as plain characters.
$findRange = $d.Content
$found = $findRange.Find.Execute("[ToC]", $false, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '[ToC]' placeholder paragraph"
}

$tocIndex = $findRange.Paragraphs.Item(1).Index

# Insert a new, empty paragraph right before the "[ToC]" paragraph, then
# fill it in with its own text/style so the placeholder paragraph itself
# is left untouched.
$targetPara = $d.Paragraphs.Item($tocIndex)
$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($tocIndex)
$newPara.Range.Text = "Table of Contents"
$newPara.Style = "Heading2"
